$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 182 ("Provincia de Limarí" /
# 44615 record), shifting rows 182-194 down to 183-195. This reproduces the
# weekly refresh: a new "Provincia de Talca" observation (44714) is added to
# the front of this block of records.
$ws.Rows.Item(182).Insert()

$ws.Cells.Item(182, 1).Value2 = 3
$ws.Cells.Item(182, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(182, 3).Value2 = "Coquimbo"
$ws.Cells.Item(182, 4).Value2 = 44714
$ws.Cells.Item(182, 5).Value2 = 5
$ws.Cells.Item(182, 6).Value2 = 100112030
$ws.Cells.Item(182, 7).Value2 = "Poroto granado"
$ws.Cells.Item(182, 8).Value2 = "Sin especificar"
$ws.Cells.Item(182, 9).Value2 = "Primera"
$ws.Cells.Item(182, 10).Value2 = 38
$ws.Cells.Item(182, 11).Value2 = 23000
$ws.Cells.Item(182, 12).Value2 = 23000
$ws.Cells.Item(182, 13).Value2 = 23000
$ws.Cells.Item(182, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(182, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(182, 16).Value2 = 920
$ws.Cells.Item(182, 17).Value2 = 25
$ws.Cells.Item(182, 18).Value2 = "Hortaliza"
